$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("theta_estimates")
$src = $wb.Worksheets.Item("top_jobs_skill")

$src.Rows.Item(2).Copy()
$ws.Rows.Item(20).PasteSpecial(-4122) | Out-Null  # formats
$excel.CutCopyMode = $false

Write-Host "done"
